$wb = $excel.ActiveWorkbook

# Rename the first sheet "Paineis DARQ" -> "PAINEIS DARQ"
$painel = $wb.Worksheets.Item("Paineis DARQ")
$painel.Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$recolh = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$recolh.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Delete the "Desarquivamentos Pendentes" sheet entirely
$desarq = $wb.Worksheets.Item("Desarquivamentos Pendentes")
$desarq.Delete()
